$p = $ppt.ActivePresentation

# Slide 2: merge the first two title runs ("Class Work #" + "1") into a
# single run reading "Class Work #1", leaving the trailing fa-IR "1" run
# (the 3rd run) untouched.
$s2 = $p.Slides.Item(2)
$title = $s2.Shapes.Item(1)
$tr = $title.TextFrame.TextRange
$merged = $tr.Characters(1, 13)
$merged.Text = "Class Work #1"

# Slide 3: un-hide the slide (removes show="0").
$s3 = $p.Slides.Item(3)
$s3.SlideShowTransition.Hidden = 0

# Slide 4: un-hide the slide (removes show="0").
$s4 = $p.Slides.Item(4)
$s4.SlideShowTransition.Hidden = 0
